$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.255.32"
$ws.Range("E2").Value = "  +14.55%  "
$ws.Range("D3").Value = "1.677.18"
$ws.Range("E3").Value = "  +8.87%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'307.58"
$ws.Range("E5").Value = "  +9.20%  "
$ws.Range("D6").Value = "'0.9954"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("D7").Value = "'0.3727"
$ws.Range("E7").Value = "  +2.89%  "
$ws.Range("D8").Value = "'0.3433"
$ws.Range("E8").Value = "  +8.10%  "
$ws.Range("D9").Value = "'48.20"
$ws.Range("E9").Value = "  +17.69%  "
$ws.Range("D10").Value = "'1.189"
$ws.Range("E10").Value = "  +8.21%  "
$ws.Range("D11").Value = "'0.07289"
$ws.Range("E11").Value = "  +7.20%  "
$ws.Range("D12").Value = "'0.9974"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'20.63"
$ws.Range("E13").Value = "  +10.12%  "
$ws.Range("D14").Value = "'6.087"
$ws.Range("E14").Value = "  +7.58%  "
$ws.Range("D15").Value = "'6.765"
$ws.Range("E15").Value = "  +6.61%  "
$ws.Range("D16").Value = "1.672.21"
$ws.Range("E16").Value = "  +8.92%  "
$ws.Range("D17").Value = "'0.00001110"
$ws.Range("E17").Value = "  +6.26%  "
$ws.Range("D18").Value = "'0.9958"
$ws.Range("E18").Value = "  +3.88%  "
$ws.Range("D19").Value = "'0.06711"
$ws.Range("E19").Value = "  +11.18%  "
$ws.Range("D20").Value = "'81.95"
$ws.Range("E20").Value = "  +13.85%  "
$ws.Range("D21").Value = "'16.49"
$ws.Range("E21").Value = "  +9.67%  "
$ws.Range("D22").Value = "'6.154"
$ws.Range("E22").Value = "  +8.39%  "
$ws.Range("E23").Value = "  +6.02%  "
$ws.Range("D24").Value = "24.213.89"
$ws.Range("E24").Value = "  +13.93%  "
$ws.Range("D25").Value = "'2.404"
$ws.Range("E25").Value = "  +4.06%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.677"
$ws.Range("E26").Value = "  +21.47%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'3.354"
$ws.Range("E27").Value = "  -9.43%  "
$ws.Range("D28").Value = "'151.89"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("D29").Value = "'19.55"
$ws.Range("E29").Value = "  +10.75%  "
$ws.Range("D30").Value = "1.852.12"
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("D31").Value = "'127.44"
$ws.Range("E31").Value = "  +8.01%  "
$ws.Range("D32").Value = "'6.361"
$ws.Range("E32").Value = "  +22.53%  "
$ws.Range("D33").Value = "'4.021"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "'0.9892"
$ws.Range("E34").Value = "  +16.23%  "
$ws.Range("D35").Value = "'1.740"
$ws.Range("E35").Value = "  +15.89%  "
$ws.Range("D36").Value = "'0.08436"
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("D37").Value = "'12.49"
$ws.Range("E37").Value = "  +16.49%  "
$ws.Range("D38").Value = "'5.373"
$ws.Range("E38").Value = "  +8.38%  "
$ws.Range("D39").Value = "'8.942"
$ws.Range("E39").Value = "  +16.26%  "
$ws.Range("D40").Value = "'0.06414"
$ws.Range("E40").Value = "  +9.41%  "
$ws.Range("D41").Value = "'1.296"
$ws.Range("E41").Value = "  +6.94%  "
$ws.Range("D42").Value = "'0.02352"
$ws.Range("E42").Value = "  +12.12%  "
$ws.Range("D43").Value = "'0.2120"
$ws.Range("E43").Value = "  +10.98%  "
$ws.Range("D44").Value = "'0.6130"
$ws.Range("E44").Value = "  +12.41%  "
$ws.Range("D45").Value = "'0.9951"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").Value = "'13.23"
$ws.Range("E46").Value = "  +6.35%  "
$ws.Range("D47").Value = "'3.799"
$ws.Range("E47").Value = "  +6.49%  "
$ws.Range("D48").Value = "'0.5976"
$ws.Range("E48").Value = "  +9.95%  "
$ws.Range("D49").Value = "'127.20"
$ws.Range("E49").Value = "  +4.70%  "
$ws.Range("D50").Value = "'2.026"
$ws.Range("E50").Value = "  +8.42%  "
$ws.Range("D51").Value = "'0.07143"
$ws.Range("E51").Value = "  +7.81%  "
